$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: update the wording/details of the upcoming presentation
$ws.Range("I29").Value = "Presentation on 27/03 at 11:30 in room A118"

# Row 17: progress 30% -> 40%, add a comment about current browser support
$ws.Range("E17").Value = 0.4
$ws.Range("I17").Value = "Currently works for chrome"

# Row 18: progress 50% -> 90%, mark "Done" column (G) with an x, add a comment
$ws.Range("E18").Value = 0.9
$ws.Range("G18").Value = "x"
$ws.Range("I18").Value = "Idk, It should be fine? Might not be 100% bug free, so I'm putting it on finished, but not 100%"

# Row 19: mark "Pending" column (D) with an x instead of a numeric progress value
$ws.Range("D19").Value = "x"
$ws.Range("E19").ClearContents()

# Move the active selection to match where the author left off editing
$ws.Range("E16").Select()
